{"js": "// Office.js (Word JavaScript API) edit script.\n// Rewrites the \"Quantum Entanglement\" article into the \"Chemistry\" article,\n// renames the author, changes the e-mail address, and appends a trailing\n// empty paragraph \u2014 matching the supplied OOXML diff.\n\n// Helper: find a unique run of text in the body and replace its contents\n// in place (keeps the original run formatting because Range.insertText\n// with Replace re-uses the formatting of the range being replaced).\nasync function replaceText(context, findText, newText) {\n  const results = context.document.body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n  return results.items[0];\n}\n\n// 1. Title\nawait replaceText(context, \"Quantum Entanglement: Unveiling the Enigma\",\n  \"Embracing the Fascinating World of Chemistry: Unveiling the Secrets of Matter\");\n\n// 2. Author name: \"Fiona Campbell\" -> \"Dr. Randall Ernest\" (built out of\n// three runs: \"Dr\" + \".\" + \" Randall Ernest\").\n{\n  const results = context.document.body.search(\"Fiona Campbell\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const nameRange = results.items[0];\n  nameRange.insertText(\"Dr\", Word.InsertLocation.replace);\n  await context.sync();\n  const dotRange = nameRange.insertText(\".\", Word.InsertLocation.after);\n  await context.sync();\n  dotRange.insertText(\" Randall Ernest\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// 3. E-mail local part + domain: \"fiona\" + \".\" + \"campbell@metaphysic-studies\"\n// collapse into a single run \"randallerest@paddingtonmail\" (the following\n// \".\" and \"org\" runs are left untouched, producing\n// \"randallerest@paddingtonmail.org\").\nawait replaceText(context, \"fiona.campbell@metaphysic-studies\", \"randallerest@paddingtonmail\");\n\n// 4. Body paragraph sentences (quantum entanglement -> chemistry).\nawait replaceText(\n  context,\n  \"In the realm of quantum mechanics, where the boundaries of conventional physics blur, lies a profound and perplexing phenomenon known as quantum entanglement\",\n  \"In the vast tapestry of scientific disciplines, chemistry stands out as a fascinating and intricate realm, holding the key to understanding the fundamental building blocks of matter\"\n);\n\nawait replaceText(\n  context,\n  \" This intricate dance of interconnectedness between subatomic particles, regardless of their physical separation, has captivated the scientific community for decades\",\n  \" Step into the enthralling world of chemistry, where curiosity kindles the flame of discovery, and the enigmatic dance of atoms and molecules unravels the mysteries of the universe\"\n);\n\nawait replaceText(\n  context,\n  \" Quantum entanglement transcends the limits of space and time, challenging our understanding of reality and opening up new avenues of exploration that could revolutionize various fields\",\n  \" Like an alchemical symphony, chemistry weaves together elements, transforming them into substances that shape our lives, from the air we breathe to the food we consume\"\n);\n\nawait replaceText(\n  context,\n  \" From the fabric of the universe to the promise of secure communication and computation, quantum entanglement holds the key to unlocking some of the most profound mysteries that govern our existence\",\n  \" Through the intricate interplay of chemical reactions, we delve into the quantum dance of particles, unlocking the secrets of matter at its core\"\n);\n\nawait replaceText(\n  context,\n  \"The intertwining of quantum particles in entanglement defies our classical notions of causality and locality\",\n  \"As we embark on this captivating journey, chemistry unveils the secrets of the universe, inviting us to explore intricate reactions and extraordinary transformations\"\n);\n\nawait replaceText(\n  context,\n  \" Measurements performed on one entangled particle instantaneously affect the properties of its distant counterpart, even if they are lightyears apart\",\n  \" From the formation of stars in the distant galaxies to the chemical interactions that sustain life on Earth, chemistry plays a pivotal role, shaping the fabric of our existence\"\n);\n\nawait replaceText(\n  context,\n  \" This nonlocal connection raises fundamental questions about the nature of reality, forcing us to confront the limits of our understanding\",\n  \" While the complexities of the universe may seem daunting, chemistry provides us with a powerful lens through which we can comprehend the interconnectedness of all matter\"\n);\n\nawait replaceText(\n  context,\n  \" The paradoxes and counterintuitive implications of quantum entanglement have sparked intense debate among physicists, leading to various interpretations and theories attempting to unravel its enigmatic essence\",\n  \" In this exploration of chemistry, we unravel the mysteries of chemical bonds, explore the properties of diverse substances, and delve into the realm of chemical reactions, painting a vivid portrait of the world around us\"\n);\n\nawait replaceText(\n  context,\n  \"As we delve deeper into the mysteries of quantum entanglement, potential applications emerge, promising to reshape our technological landscape\",\n  \"The world of chemistry promises an exhilarating adventure, filled with wonder, discovery, and practical applications\"\n);\n\nawait replaceText(\n  context,\n  \" Quantum cryptography, harnessing the inherent randomness of entangled particles, offers unbreakable encryption methods\",\n  \" By understanding the underlying principles that govern chemical reactions, we gain valuable insights into diverse phenomena, ranging from the intricate workings of pharmaceuticals to the intricacies of industrial processes\"\n);\n\nawait replaceText(\n  context,\n  \" Quantum computing, exploiting the superposition and entanglement properties, holds the promise of exponential leaps in computational power\",\n  \" As we delve deeper into this fascinating field, we are empowered to contribute meaningfully to addressing global challenges, such as developing sustainable energy sources, finding cures for diseases, and creating innovative materials that shape the future of technology\"\n);\n\n// The next three runs (\" Furthermore...universe\", \".\", \" The study of\n// quantum \") collapse into a single run \" This pursuit of chemical \".\nawait replaceText(\n  context,\n  \" Furthermore, the intricate dance of entangled particles may play a crucial role in unraveling the enigmas of gravitational interactions and the fundamental forces that shape our universe. The study of quantum \",\n  \" This pursuit of chemical \"\n);\n\nawait replaceText(\n  context,\n  \"entanglement represents a scientific odyssey, pushing the boundaries of human knowledge and opening up new horizons of discovery\",\n  \"knowledge is not merely an academic exercise but a testament to our insatiable curiosity and unwavering commitment to understanding the fundamental fabric of the universe\"\n);\n\n// 5. Summary paragraph.\nawait replaceText(\n  context,\n  \"Quantum entanglement, an enigmatic phenomenon in quantum mechanics, challenges our understanding of reality with its nonlocal correlations and profound implications\",\n  \"Chemistry, a captivating discipline, unveils the secrets of matter, unraveling the intricate tapestry of elements, compounds, and reactions\"\n);\n\nawait replaceText(\n  context,\n  \" This intricate interconnectedness between subatomic particles holds the potential to revolutionize diverse fields, from secure communication and computing to fundamental physics\",\n  \" Through the study of chemistry, we delve into the quantum dance of atoms and molecules, exploring the fundamental principles that govern chemical interactions\"\n);\n\n// The final sentence is replaced and then two brand-new runs (a \".\"\n// followed by a new closing sentence) are appended before the pre-existing\n// trailing \".\" run.\n{\n  const find = \" As we unravel the mysteries of quantum entanglement, we embark on a captivating journey into the heart of quantum mechanics, pushing the boundaries of knowledge and innovation\";\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  const r = results.items[0];\n  r.insertText(\n    \" Chemistry provides a powerful lens through which we comprehend the interconnectedness of all matter, empowering us to address global challenges and shape the future of technology\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n  const dotRange = r.insertText(\".\", Word.InsertLocation.after);\n  await context.sync();\n  dotRange.insertText(\n    \" Embracing the fascinating world of chemistry, we embark on an extraordinary adventure filled with wonder, discovery, and practical applications, ultimately enriching our understanding of the universe and our place within it\",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n}\n\n// 6. Append a new trailing empty paragraph at the very end of the body.\ncontext.document.body.insertParagraph(\"\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Word COM interop script.\n# Rewrites the \"Quantum Entanglement\" article into the \"Chemistry\" article,\n# renames the author, changes the e-mail address, and appends a trailing\n# empty paragraph \u2014 matching the supplied OOXML diff.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $newText) {\n    $r = $d.Content\n    $r.Find.MatchCase = $true\n    $null = $r.Find.Execute($findText)\n    $r.Text = $newText\n    return $r\n}\n\n# 1. Title\nReplace-Text \"Quantum Entanglement: Unveiling the Enigma\" \"Embracing the Fascinating World of Chemistry: Unveiling the Secrets of Matter\" | Out-Null\n\n# 2. Author name: \"Fiona Campbell\" -> \"Dr. Randall Ernest\" (built out of\n# three runs: \"Dr\" + \".\" + \" Randall Ernest\").\n$r = $d.Content\n$r.Find.MatchCase = $true\n$null = $r.Find.Execute(\"Fiona Campbell\")\n$r.Text = \"Dr\"\n$r.Collapse(0)\n$r.InsertAfter(\".\")\n$r.Collapse(0)\n$r.InsertAfter(\" Randall Ernest\")\n\n# 3. E-mail local part + domain: \"fiona\" + \".\" + \"campbell@metaphysic-studies\"\n# collapse into a single run \"randallerest@paddingtonmail\" (the following\n# \".\" and \"org\" runs are left untouched, producing\n# \"randallerest@paddingtonmail.org\").\nReplace-Text \"fiona.campbell@metaphysic-studies\" \"randallerest@paddingtonmail\" | Out-Null\n\n# 4. Body paragraph sentences (quantum entanglement -> chemistry).\nReplace-Text \"In the realm of quantum mechanics, where the boundaries of conventional physics blur, lies a profound and perplexing phenomenon known as quantum entanglement\" \"In the vast tapestry of scientific disciplines, chemistry stands out as a fascinating and intricate realm, holding the key to understanding the fundamental building blocks of matter\" | Out-Null\n\nReplace-Text \" This intricate dance of interconnectedness between subatomic particles, regardless of their physical separation, has captivated the scientific community for decades\" \" Step into the enthralling world of chemistry, where curiosity kindles the flame of discovery, and the enigmatic dance of atoms and molecules unravels the mysteries of the universe\" | Out-Null\n\nReplace-Text \" Quantum entanglement transcends the limits of space and time, challenging our understanding of reality and opening up new avenues of exploration that could revolutionize various fields\" \" Like an alchemical symphony, chemistry weaves together elements, transforming them into substances that shape our lives, from the air we breathe to the food we consume\" | Out-Null\n\nReplace-Text \" From the fabric of the universe to the promise of secure communication and computation, quantum entanglement holds the key to unlocking some of the most profound mysteries that govern our existence\" \" Through the intricate interplay of chemical reactions, we delve into the quantum dance of particles, unlocking the secrets of matter at its core\" | Out-Null\n\nReplace-Text \"The intertwining of quantum particles in entanglement defies our classical notions of causality and locality\" \"As we embark on this captivating journey, chemistry unveils the secrets of the universe, inviting us to explore intricate reactions and extraordinary transformations\" | Out-Null\n\nReplace-Text \" Measurements performed on one entangled particle instantaneously affect the properties of its distant counterpart, even if they are lightyears apart\" \" From the formation of stars in the distant galaxies to the chemical interactions that sustain life on Earth, chemistry plays a pivotal role, shaping the fabric of our existence\" | Out-Null\n\nReplace-Text \" This nonlocal connection raises fundamental questions about the nature of reality, forcing us to confront the limits of our understanding\" \" While the complexities of the universe may seem daunting, chemistry provides us with a powerful lens through which we can comprehend the interconnectedness of all matter\" | Out-Null\n\nReplace-Text \" The paradoxes and counterintuitive implications of quantum entanglement have sparked intense debate among physicists, leading to various interpretations and theories attempting to unravel its enigmatic essence\" \" In this exploration of chemistry, we unravel the mysteries of chemical bonds, explore the properties of diverse substances, and delve into the realm of chemical reactions, painting a vivid portrait of the world around us\" | Out-Null\n\nReplace-Text \"As we delve deeper into the mysteries of quantum entanglement, potential applications emerge, promising to reshape our technological landscape\" \"The world of chemistry promises an exhilarating adventure, filled with wonder, discovery, and practical applications\" | Out-Null\n\nReplace-Text \" Quantum cryptography, harnessing the inherent randomness of entangled particles, offers unbreakable encryption methods\" \" By understanding the underlying principles that govern chemical reactions, we gain valuable insights into diverse phenomena, ranging from the intricate workings of pharmaceuticals to the intricacies of industrial processes\" | Out-Null\n\nReplace-Text \" Quantum computing, exploiting the superposition and entanglement properties, holds the promise of exponential leaps in computational power\" \" As we delve deeper into this fascinating field, we are empowered to contribute meaningfully to addressing global challenges, such as developing sustainable energy sources, finding cures for diseases, and creating innovative materials that shape the future of technology\" | Out-Null\n\n# The next three runs (\" Furthermore...universe\", \".\", \" The study of\n# quantum \") collapse into a single run \" This pursuit of chemical \".\nReplace-Text \" Furthermore, the intricate dance of entangled particles may play a crucial role in unraveling the enigmas of gravitational interactions and the fundamental forces that shape our universe. The study of quantum \" \" This pursuit of chemical \" | Out-Null\n\nReplace-Text \"entanglement represents a scientific odyssey, pushing the boundaries of human knowledge and opening up new horizons of discovery\" \"knowledge is not merely an academic exercise but a testament to our insatiable curiosity and unwavering commitment to understanding the fundamental fabric of the universe\" | Out-Null\n\n# 5. Summary paragraph.\nReplace-Text \"Quantum entanglement, an enigmatic phenomenon in quantum mechanics, challenges our understanding of reality with its nonlocal correlations and profound implications\" \"Chemistry, a captivating discipline, unveils the secrets of matter, unraveling the intricate tapestry of elements, compounds, and reactions\" | Out-Null\n\nReplace-Text \" This intricate interconnectedness between subatomic particles holds the potential to revolutionize diverse fields, from secure communication and computing to fundamental physics\" \" Through the study of chemistry, we delve into the quantum dance of atoms and molecules, exploring the fundamental principles that govern chemical interactions\" | Out-Null\n\n# The final sentence is replaced and then two brand-new runs (a \".\"\n# followed by a new closing sentence) are appended before the pre-existing\n# trailing \".\" run.\n$r = $d.Content\n$r.Find.MatchCase = $true\n$null = $r.Find.Execute(\" As we unravel the mysteries of quantum entanglement, we embark on a captivating journey into the heart of quantum mechanics, pushing the boundaries of knowledge and innovation\")\n$r.Text = \" Chemistry provides a powerful lens through which we comprehend the interconnectedness of all matter, empowering us to address global challenges and shape the future of technology\"\n$r.Collapse(0)\n$r.InsertAfter(\".\")\n$r.Collapse(0)\n$r.InsertAfter(\" Embracing the fascinating world of chemistry, we embark on an extraordinary adventure filled with wonder, discovery, and practical applications, ultimately enriching our understanding of the universe and our place within it\")\n\n# 6. Append a new trailing empty paragraph at the very end of the document.\n$end = $d.Content\n$end.Collapse(0)\n$end.InsertParagraphAfter()\n"}
